$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '62.759.16'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +3.11%  '
$c.Style = 'Normal'

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.446.72'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +1.94%  '
$c.Style = 'Normal'

$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -0.19%  '
$c.Style = 'Normal'

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '576.03'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  +1.89%  '
$c.Style = 'Normal'

$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  +2.60%  '
$c.Style = 'Normal'

$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c.Style = 'Normal'

$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'

$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '2.445.08'
$c.Style = 'Normal'
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +1.56%  '
$c.Style = 'Normal'

$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +2.55%  '
$c.Style = 'Normal'

$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  +2.48%  '
$c.Style = 'Normal'

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '5.27'
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +1.45%  '
$c.Style = 'Normal'

$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  +1.89%  '
$c.Style = 'Normal'

$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '28.27'
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +7.12%  '
$c.Style = 'Normal'

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '0.0000179'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +5.13%  '
$c.Style = 'Normal'

$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '2.887.66'
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  +1.70%  '
$c.Style = 'Normal'

$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '62.623.84'
$c.Style = 'Normal'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  +3.19%  '
$c.Style = 'Normal'

$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '2.444.56'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c.Style = 'Normal'

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '7.95'
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -1.12%  '
$c.Style = 'Normal'

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '11.00'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  +2.78%  '
$c.Style = 'Normal'

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '330.18'
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +1.82%  '
$c.Style = 'Normal'

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.Style = 'Normal'

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '66.23'
$c.Style = 'Normal'

$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '644.92'
$c.Style = 'Normal'
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +9.66%  '
$c.Style = 'Normal'

$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  +17.73%  '
$c.Style = 'Normal'

$c = $ws.Range("B28")
$c.NumberFormat = '@'
$c.Value = 'Aptos'
$c.Style = 'Normal'
$c = $ws.Range("C28")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c.Style = 'Normal'
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '8.49'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  +3.30%  '
$c.Style = 'Normal'

$c = $ws.Range("B29")
$c.NumberFormat = '@'
$c.Value = 'BabyDogeCoin'
$c.Style = 'Normal'
$c = $ws.Range("C29")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c.Style = 'Normal'
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '0.0₆0531'
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +90.71%  '
$c.Style = 'Normal'

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '0.0₃0990'
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  +4.98%  '
$c.Style = 'Normal'

$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '2.567.29'
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +2.31%  '
$c.Style = 'Normal'

$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '8.21'
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  +2.21%  '
$c.Style = 'Normal'

$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  +6.75%  '
$c.Style = 'Normal'

$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +2.70%  '
$c.Style = 'Normal'

$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  +4.24%  '
$c.Style = 'Normal'

$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +1.30%  '
$c.Style = 'Normal'

$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.Style = 'Normal'

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '4.76'
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +3.14%  '
$c.Style = 'Normal'

$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  +5.74%  '
$c.Style = 'Normal'

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '153.51'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  +1.08%  '
$c.Style = 'Normal'

$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +0.76%  '
$c.Style = 'Normal'

$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '18.73'
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  +2.23%  '
$c.Style = 'Normal'

$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '2.72'
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  +7.94%  '
$c.Style = 'Normal'

$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '1.76'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  +4.14%  '
$c.Style = 'Normal'

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '42.47'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +2.09%  '
$c.Style = 'Normal'

$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.Style = 'Normal'

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '14.94'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  +27.42%  '
$c.Style = 'Normal'

$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '145.27'
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  +2.39%  '
$c.Style = 'Normal'

$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +3.33%  '
$c.Style = 'Normal'

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '20.62'
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  +5.81%  '
$c.Style = 'Normal'

$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  +2.31%  '
$c.Style = 'Normal'
